# Bump the "Förändrad" (Changed) date in column C by one day for all data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
